$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships with protection enabled; temporarily unprotect so the
# cells below (which are locked, like the rest of the sheet) can be edited.
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure note.
$ws.Range("A12").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-18 for illustrative purposes only and are subject to change."

# Refresh the Weight (D) and Percent Change (E) columns with the latest
# model-holdings figures.
$ws.Range("D2").Value = 0.1778300754914373
$ws.Range("E2").Value = -0.002710027100271128

$ws.Range("D3").Value = 0.1776674260321464
$ws.Range("E3").Value = -0.002949852507374562

$ws.Range("D4").Value = 0.2243608729657523
$ws.Range("E4").Value = -0.005020920502091908

$ws.Range("D5").Value = 0.08002353397114678
$ws.Range("E5").Value = -0.00101626016260159

$ws.Range("D6").Value = 0.07999140568289179

$ws.Range("D7").Value = 0.1203686319473658
$ws.Range("E7").Value = -0.0009842519685039353

$ws.Range("D8").Value = 0.1397580539092597
$ws.Range("E8").Value = -0.001666666666666594

$ws.Range("E9").Value = -0.002565243015360208

# Restore sheet protection to match the shipped state.
$ws.Protect()
